$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fill in the "Beat Vegas?" (column G) results for the existing
#    games played on 2021-01-20 (rows 99-107), which were left blank
#    until the outcomes were known.
# ------------------------------------------------------------------
$ws.Range("G99").Value  = "No"
$ws.Range("G100").Value = "Yes"
$ws.Range("G101").Value = "No"
$ws.Range("G102").Value = "No"
$ws.Range("G103").Value = "Yes"
$ws.Range("G104").Value = "Yes"
$ws.Range("G105").Value = "No"
$ws.Range("G106").Value = "No"
$ws.Range("G107").Value = "Yes"

# ------------------------------------------------------------------
# 2) Append the new games scheduled for 1/21/2021 (rows 108-110).
#    Copy the formatting from the last existing data row (107) down
#    into the new rows first so the date cells keep the same
#    yyyy-mm-dd number format / style as the rest of column A.
# ------------------------------------------------------------------
$ws.Range("A107:G107").Copy($ws.Range("A108:G110"))

$ws.Range("A108").Value = 44217
$ws.Range("B108").Value = "MIL"
$ws.Range("C108").Value = "LAL"
$ws.Range("D108").Value = -1
$ws.Range("E108").Value = 0
$ws.Range("F108").Value = -1
$ws.Range("G108").Value = ""

$ws.Range("A109").Value = 44217
$ws.Range("B109").Value = "NOP"
$ws.Range("C109").Value = "UTA"
$ws.Range("D109").Value = -7
$ws.Range("E109").Value = 7.5
$ws.Range("F109").Value = -14.5
$ws.Range("G109").Value = ""

$ws.Range("A110").Value = 44217
$ws.Range("B110").Value = "NYK"
$ws.Range("C110").Value = "GSW"
$ws.Range("D110").Value = -4.5
$ws.Range("E110").Value = 11
$ws.Range("F110").Value = -15.5
$ws.Range("G110").Value = ""

# ------------------------------------------------------------------
# 3) Restore the view state roughly matching where the user left the
#    selection after entering the new rows.
# ------------------------------------------------------------------
$ws.Range("O97").Select()
